$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select column F (numberOfPages) and delete it entirely, shifting
# documentType (G) and abc (H) one column to the left.
$col = $ws.Range("F1:F1048576")
$col.Select()
$col.Delete()
